# Updated symbol list on Mon Jan 16 09:55:21 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns with latest quotes.
# Values are written with a leading apostrophe so Excel keeps them as
# literal text (matching the workbook's existing inlineStr/text cells)
# instead of auto-converting to numbers or percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.48"
$ws.Range("E2").Value = "'1.52%"
$ws.Range("D3").Value = "'31.30"
$ws.Range("E3").Value = "'0.09%"
$ws.Range("D4").Value = "'5.145"
$ws.Range("E4").Value = "'0.94%"
$ws.Range("D5").Value = "'0.08052"
$ws.Range("E5").Value = "'9.44%"
$ws.Range("D6").Value = "'2.725"
$ws.Range("E6").Value = "'65.31%"
$ws.Range("D7").Value = "'7.853"
$ws.Range("E7").Value = "'2.22%"
$ws.Range("D8").Value = "'3.821"
$ws.Range("E8").Value = "'1.98%"
$ws.Range("D9").Value = "'0.9086"
$ws.Range("E9").Value = "'-1.20%"
$ws.Range("E10").Value = "'3.36%"
$ws.Range("D11").Value = "'0.07218"
$ws.Range("E11").Value = "'1.37%"
$ws.Range("D12").Value = "'0.08049"
$ws.Range("E12").Value = "'1.38%"
$ws.Range("D13").Value = "'0.03016"
$ws.Range("E13").Value = "'0.99%"
$ws.Range("D14").Value = "'0.09967"
$ws.Range("E14").Value = "'0.84%"
$ws.Range("D15").Value = "'0.001493"
$ws.Range("E15").Value = "'-0.05%"
$ws.Range("D16").Value = "'0.005935"
$ws.Range("E16").Value = "'-4.37%"
$ws.Range("D17").Value = "'3.503"
$ws.Range("D18").Value = "'2.241"
$ws.Range("E18").Value = "'0.55%"
$ws.Range("D20").Value = "'0.1330"
$ws.Range("E20").Value = "'-0.28%"
$ws.Range("D21").Value = "'4.598"
$ws.Range("E21").Value = "'0.92%"
$ws.Range("E22").Value = "'3.22%"
$ws.Range("D23").Value = "'0.04605"
$ws.Range("E23").Value = "'-0.41%"
$ws.Range("D24").Value = "'0.001264"
$ws.Range("E24").Value = "'3.81%"
$ws.Range("D25").Value = "'0.004447"
$ws.Range("E25").Value = "'0.63%"
$ws.Range("D26").Value = "'0.0001181"
$ws.Range("E26").Value = "'-9.06%"
$ws.Range("D27").Value = "'0.0003433"
$ws.Range("E27").Value = "'83.08%"
$ws.Range("D39").Value = "'0.01817"
$ws.Range("E39").Value = "'7.18%"
$ws.Range("D40").Value = "'0.04534"
$ws.Range("E40").Value = "'3.00%"
$ws.Range("D41").Value = "'0.007100"
$ws.Range("E41").Value = "'0.26%"
$ws.Range("D42").Value = "'0.1342"
$ws.Range("E42").Value = "'1.16%"
$ws.Range("D43").Value = "'0.002242"
$ws.Range("E43").Value = "'6.35%"
$ws.Range("D44").Value = "'0.01044"
$ws.Range("E44").Value = "'-5.15%"
$ws.Range("D45").Value = "'0.00006357"
$ws.Range("E45").Value = "'6.13%"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("E47").Value = "'-53.83%"
$ws.Range("D48").Value = "'0.006205"
$ws.Range("E48").Value = "'-43.57%"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.14%"
